# OW-535 updated the bilateral trade portfolio to match the acuo-data test branch
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IRS-Bilateral")

# Position Account ID (row 2) changed to a new ACUO identifier
$ws.Range("B2").Value = "ACUOSG8745"

# Portfolio ID (row 2) changed from p9 to p1
$ws.Range("AP2").Value = "p1"

# Reset the view back to the top-left of the sheet and select A2 only
$ws.Range("A2").Select()
